# Update "想去人数" (interested count) and "最低票价" (min ticket price) figures
# for the 展览 (Exhibitions) and 全部类型 (All types) sheets — both sheets carry
# identical event listings, so the same row/column edits are applied to each.

$wb = $excel.ActiveWorkbook

# Map of row -> column -> new value, applied identically on both sheets.
$updates = @(
    @{ Row = 2;  Col = "G"; Value = 70 },
    @{ Row = 6;  Col = "F"; Value = 40 },
    @{ Row = 7;  Col = "F"; Value = 12180 },
    @{ Row = 11; Col = "F"; Value = 424 },
    @{ Row = 12; Col = "F"; Value = 1125 },
    @{ Row = 13; Col = "F"; Value = 883 },
    @{ Row = 14; Col = "F"; Value = 13560 },
    @{ Row = 15; Col = "F"; Value = 13659 },
    @{ Row = 17; Col = "F"; Value = 159 },
    @{ Row = 21; Col = "F"; Value = 98 },
    @{ Row = 23; Col = "F"; Value = 2113 },
    @{ Row = 24; Col = "F"; Value = 195 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $addr = "$($u.Col)$($u.Row)"
        $ws.Range($addr).Value = $u.Value
    }
}
